# Auto-save via app Streamlit
# A new reservation ("Marco Sambaan") was inserted as the first data row
# (row 2) of the reservations sheet; every subsequent row (old rows 2-51)
# shifts down by one, so the TOTAL row moves from row 51 to row 52 and the
# used range grows from A1:O51 to A1:O52.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row above the current row 2, pushing all existing
# reservations (and the TOTAL row) down by one row.
$ws.Rows("2:2").Insert()

# Excel's default row-insert behaviour copies the formatting of the row
# above (the bold/centered header row here). Strip that back off so the
# new row matches the plain look of the other data rows.
$ws.Rows("2:2").ClearFormats()

# Re-apply the date formatting used by the other rows' arrival/departure
# date columns.
$ws.Range("D2:E2").NumberFormat = "YYYY-MM-DD"

# Fill in the new reservation's data.
$ws.Range("A2").Value = "Marco Sambaan"
$ws.Range("B2").Value = "Booking"
# C2 (telephone) is intentionally left blank for this reservation.
$ws.Range("D2").Value = 45518
$ws.Range("E2").Value = 45881
$ws.Range("F2").Value = 363
$ws.Range("G2").Value = 397
$ws.Range("H2").Value = 223.78
$ws.Range("I2").Value = 173.22
$ws.Range("J2").Value = 43.63
$ws.Range("K2").Value = 2024
$ws.Range("L2").Value = 8
# M2/N2/O2 (uid_ical, annee, mois) remain blank for this reservation.
